$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47: shorten the farewell line, enable wrap text, and add the "next" trigger
$ws.Range("E47").Value = "그럼 튜토리얼은 여기까지에요. 즐거운 시간이셨는지 모르겠네요."
$ws.Range("E47").WrapText = $true
$ws.Range("G47").Value = "next"

# New row 48: jigi says goodbye to the visitor
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "dialogue"
$ws.Cells.Item(48, 3).Value = "jigi"
$ws.Cells.Item(48, 5).Value = "방문해주셔서 감사합니다, 내일도 또 와주세요!"

# Move the active selection to match the edited workbook
$ws.Range("G47").Select()
